{"js": "/* Replace the arithmetic-expression text in each cell of the single\n * table (20 rows x 5 columns = 100 cells) with its updated value, in\n * document order. A couple of \"before\" strings repeat verbatim (e.g.\n * \"99-55=\" appears twice, mapped to two different replacements), so\n * cells must be addressed positionally (row, col), not via text search\n * / find-replace. REPLACEMENTS[i] = [oldText, newText] in row-major\n * (left-to-right, top-to-bottom) document order. */\nconst REPLACEMENTS = [[\"89-32=\", \"99-2=\"], [\"74-73=\", \"87-44=\"], [\"95-37=\", \"28+43=\"], [\"95-25=\", \"71+8=\"], [\"19+71=\", \"47+10=\"], [\"20+26=\", \"97-60=\"], [\"6+20=\", \"6+75=\"], [\"7+32=\", \"4+6=\"], [\"78-19=\", \"65-31=\"], [\"89-15=\", \"75+11=\"], [\"25-7=\", \"96-90=\"], [\"25+53=\", \"25+61=\"], [\"20+24=\", \"64-22=\"], [\"50-15=\", \"95-4=\"], [\"41-31=\", \"7+60=\"], [\"63-2=\", \"53-37=\"], [\"25+41=\", \"38+52=\"], [\"98-30=\", \"50+19=\"], [\"72-25=\", \"51+46=\"], [\"75-10=\", \"85-42=\"], [\"73-13=\", \"39+3=\"], [\"54-8=\", \"68+3=\"], [\"33-0=\", \"80+6=\"], [\"70+10=\", \"71-41=\"], [\"79-28=\", \"2+16=\"], [\"13+21=\", \"33-33=\"], [\"30-17=\", \"12+68=\"], [\"68+22=\", \"34-23=\"], [\"86-3=\", \"26+63=\"], [\"0+61=\", \"35+24=\"], [\"17+50=\", \"33-25=\"], [\"87-45=\", \"45-8=\"], [\"93-8=\", \"70+20=\"], [\"55+9=\", \"68+0=\"], [\"85-21=\", \"20-13=\"], [\"55+27=\", \"39+58=\"], [\"29+36=\", \"97-20=\"], [\"79-72=\", \"30-12=\"], [\"82-63=\", \"75+8=\"], [\"45+8=\", \"9+66=\"], [\"7+81=\", \"16+39=\"], [\"7-4=\", \"3+41=\"], [\"91-34=\", \"4+95=\"], [\"53-43=\", \"73+18=\"], [\"49-23=\", \"63+26=\"], [\"59+22=\", \"8+67=\"], [\"47-13=\", \"49+25=\"], [\"92-23=\", \"28+52=\"], [\"62-38=\", \"54+7=\"], [\"58+27=\", \"52+45=\"], [\"88-70=\", \"30-10=\"], [\"86-1=\", \"71-34=\"], [\"81+5=\", \"80-22=\"], [\"56-41=\", \"22+55=\"], [\"83-65=\", \"54-40=\"], [\"61-6=\", \"15-6=\"], [\"75-44=\", \"14+41=\"], [\"40-16=\", \"88-54=\"], [\"58-5=\", \"93-4=\"], [\"85+11=\", \"50+0=\"], [\"99-55=\", \"92-47=\"], [\"42+20=\", \"6+48=\"], [\"13+49=\", \"19-9=\"], [\"2+66=\", \"90-4=\"], [\"41+11=\", \"82+3=\"], [\"10+74=\", \"15+45=\"], [\"86-20=\", \"58-14=\"], [\"26+11=\", \"30-24=\"], [\"16+12=\", \"99-76=\"], [\"79+10=\", \"10+5=\"], [\"95+3=\", \"26-19=\"], [\"33+22=\", \"18+2=\"], [\"28+28=\", \"89+0=\"], [\"37+32=\", \"99-64=\"], [\"99-55=\", \"7+8=\"], [\"71+0=\", \"31-24=\"], [\"62-32=\", \"78-27=\"], [\"25+45=\", \"12+79=\"], [\"1+47=\", \"71-5=\"], [\"8-4=\", \"38+15=\"], [\"61-61=\", \"43+45=\"], [\"61+27=\", \"92-40=\"], [\"84-14=\", \"7+51=\"], [\"42+34=\", \"35-5=\"], [\"18+27=\", \"14+11=\"], [\"77+11=\", \"37+47=\"], [\"62-3=\", \"25+64=\"], [\"7+48=\", \"13+61=\"], [\"27+47=\", \"16+24=\"], [\"36+52=\", \"68-1=\"], [\"54-42=\", \"51+24=\"], [\"88-65=\", \"24+10=\"], [\"94-34=\", \"46-25=\"], [\"42-19=\", \"97-75=\"], [\"61+31=\", \"57-26=\"], [\"78+17=\", \"28+14=\"], [\"83-18=\", \"38+51=\"], [\"89-9=\", \"6+64=\"], [\"12+69=\", \"24+48=\"], [\"31+24=\", \"79-20=\"]];\n\nconst table = context.document.body.tables.getFirstOrNullObject();\ntable.load(\"rowCount\");\nawait context.sync();\n\nif (table.isNullObject) {\n  throw new Error(\"Expected a table in the document body, found none.\");\n}\n\nconst firstRow = table.rows.getFirstOrNullObject();\nfirstRow.load(\"cellCount\");\nawait context.sync();\nconst columns = firstRow.cellCount || 5;\n\nlet i = 0;\nfor (let r = 0; r < table.rowCount && i < REPLACEMENTS.length; r++) {\n  for (let c = 0; c < columns && i < REPLACEMENTS.length; c++) {\n    const [, after] = REPLACEMENTS[i];\n    table.getCell(r, c).value = after;\n    i++;\n  }\n}\n\nawait context.sync();\n", "ps1": "# Replace the arithmetic-expression text in each cell of the single\n# table (20 rows x 5 columns = 100 cells) with its updated value, in\n# document order. A couple of \"before\" strings repeat (e.g. \"99-55=\"\n# appears twice, mapping to two different replacements), so cells are\n# addressed positionally (row, col) rather than via text search.\n$after = @(\n    \"99-2=\", \"87-44=\", \"28+43=\", \"71+8=\", \"47+10=\", \"97-60=\", \"6+75=\", \"4+6=\",\n    \"65-31=\", \"75+11=\", \"96-90=\", \"25+61=\", \"64-22=\", \"95-4=\", \"7+60=\", \"53-37=\",\n    \"38+52=\", \"50+19=\", \"51+46=\", \"85-42=\", \"39+3=\", \"68+3=\", \"80+6=\", \"71-41=\",\n    \"2+16=\", \"33-33=\", \"12+68=\", \"34-23=\", \"26+63=\", \"35+24=\", \"33-25=\", \"45-8=\",\n    \"70+20=\", \"68+0=\", \"20-13=\", \"39+58=\", \"97-20=\", \"30-12=\", \"75+8=\", \"9+66=\",\n    \"16+39=\", \"3+41=\", \"4+95=\", \"73+18=\", \"63+26=\", \"8+67=\", \"49+25=\", \"28+52=\",\n    \"54+7=\", \"52+45=\", \"30-10=\", \"71-34=\", \"80-22=\", \"22+55=\", \"54-40=\", \"15-6=\",\n    \"14+41=\", \"88-54=\", \"93-4=\", \"50+0=\", \"92-47=\", \"6+48=\", \"19-9=\", \"90-4=\",\n    \"82+3=\", \"15+45=\", \"58-14=\", \"30-24=\", \"99-76=\", \"10+5=\", \"26-19=\", \"18+2=\",\n    \"89+0=\", \"99-64=\", \"7+8=\", \"31-24=\", \"78-27=\", \"12+79=\", \"71-5=\", \"38+15=\",\n    \"43+45=\", \"92-40=\", \"7+51=\", \"35-5=\", \"14+11=\", \"37+47=\", \"25+64=\", \"13+61=\",\n    \"16+24=\", \"68-1=\", \"51+24=\", \"24+10=\", \"46-25=\", \"97-75=\", \"57-26=\", \"28+14=\",\n    \"38+51=\", \"6+64=\", \"24+48=\", \"79-20=\"\n)\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n$cols = $t.Columns.Count\n$i = 0\nfor ($r = 1; $r -le $t.Rows.Count; $r++) {\n    for ($c = 1; $c -le $cols; $c++) {\n        if ($i -ge $after.Length) { break }\n        $t.Cell($r, $c).Range.Text = $after[$i]\n        $i++\n    }\n}\n\n"}
